$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 114, shifting existing rows 114:190 down to 115:191
# (dimension grows from A1:R190 to A1:R191).
$ws.Rows(114).Insert()

# Populate the newly inserted row 114 with the new data record.
$ws.Cells.Item(114, 1).Value = 10
$ws.Cells.Item(114, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(114, 3).Value = "La Araucanía"
$ws.Cells.Item(114, 4).Value = 44447
$ws.Cells.Item(114, 5).Value = 9
$ws.Cells.Item(114, 6).Value = 100112037
$ws.Cells.Item(114, 7).Value = "Cebollín"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 55
$ws.Cells.Item(114, 11).Value = 5000
$ws.Cells.Item(114, 12).Value = 5000
$ws.Cells.Item(114, 13).Value = 5000
$ws.Cells.Item(114, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(114, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(114, 16).Value = 417
$ws.Cells.Item(114, 17).Value = 12
$ws.Cells.Item(114, 18).Value = "Hortaliza"
